$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9465
$ws.Range("D2").Value = 8368
$ws.Range("E2").Value = 0.8840993132593766
$ws.Range("F2").Value = 0.8828866849546317
$ws.Range("G2").Value = 0.0958304262959481
$ws.Range("H2").Value = 0.08460740739021877
$ws.Range("I2").Value = 40617149.71070025
$ws.Range("J2").Value = 14113742.80200613
$ws.Range("L2").Value = 14113742.80200613
$ws.Range("M2").Value = 54730892.51270638
$ws.Range("N2").Value = 800741935.3172001
$ws.Range("O2").Value = 783042128.3132001
$ws.Range("P2").Value = 0.01762583196846711
$ws.Range("Q2").Value = 0.01802424453510492

$ws.Range("C3").Value = 9642
$ws.Range("D3").Value = 8557
$ws.Range("E3").Value = 0.8874714789462768
$ws.Range("F3").Value = 0.8852679495137595
$ws.Range("G3").Value = 0.09415192489269762
$ws.Range("H3").Value = 0.08334968149253193
$ws.Range("I3").Value = 42347526.09196078
$ws.Range("J3").Value = 14709013.98571959
$ws.Range("L3").Value = 14709013.98571959
$ws.Range("M3").Value = 57056540.07768038
$ws.Range("N3").Value = 836091164.5578281
$ws.Range("O3").Value = 818610988.533758
$ws.Range("P3").Value = 0.01759259589054328
$ws.Range("Q3").Value = 0.0179682586622315

$ws.Range("C4").Value = 9836
$ws.Range("D4").Value = 8714
$ws.Range("E4").Value = 0.8859292395282635
$ws.Range("F4").Value = 0.8839521201054981
$ws.Range("G4").Value = 0.09305289739007716
$ws.Range("H4").Value = 0.08225430592991807
$ws.Range("I4").Value = 44259719.44997451
$ws.Range("J4").Value = 15345188.53406264
$ws.Range("L4").Value = 15345188.53406264
$ws.Range("M4").Value = 59604907.98403715
$ws.Range("N4").Value = 874653871.4285319
$ws.Range("O4").Value = 857204923.4225781
$ws.Range("P4").Value = 0.01754429841944225
$ws.Range("Q4").Value = 0.017901423702507

$ws.Range("C5").Value = 10027
$ws.Range("D5").Value = 8908
$ws.Range("E5").Value = 0.8884013164455968
$ws.Range("F5").Value = 0.8860155162124528
$ws.Range("G5").Value = 0.09171498176794257
$ws.Range("H5").Value = 0.08126089691553934
$ws.Range("I5").Value = 46251536.29016398
$ws.Range("J5").Value = 16008023.2459684
$ws.Range("L5").Value = 16008023.2459684
$ws.Range("M5").Value = 62259559.53613237
$ws.Range("N5").Value = 913201104.5186434
$ws.Range("O5").Value = 895714999.0556703
$ws.Range("P5").Value = 0.017529570613481
$ws.Range("Q5").Value = 0.017871782054387

$ws.Range("C6").Value = 10222
$ws.Range("D6").Value = 9062
$ws.Range("E6").Value = 0.8865192721580903
$ws.Range("F6").Value = 0.8837526818802418
$ws.Range("G6").Value = 0.09064143435825292
$ws.Range("H6").Value = 0.08010461070357792
$ws.Range("I6").Value = 48302938.50549269
$ws.Range("J6").Value = 16662679.69380575
$ws.Range("L6").Value = 16662679.69380575
$ws.Range("M6").Value = 64965618.19929844
$ws.Range("N6").Value = 954323256.6988841
$ws.Range("O6").Value = 936731430.8257025
$ws.Range("P6").Value = 0.01746020499536385
$ws.Range("Q6").Value = 0.01778810782415838
